$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 93. This pushes the old blank row 93
# down to row 95, and the old "B2F" row 94 down to row 96 - exactly
# like Excel's normal "insert rows above" behaviour.
$ws.Range("A93:H94").Insert()

# --- New row 93: "Affine Texture" ---
$ws.Range("A93").Value = "Affine Texture"
$ws.Range("B93").Value = 265
$ws.Range("C93").Formula = "=B93/30"
$ws.Range("D93").Formula = "=B93/B`$85"
$ws.Range("E93").Value = 252
$ws.Range("F93").Value = 909
$ws.Range("G93").Formula = "=E93*C93"
$ws.Range("H93").Formula = "=F93*C93"

# --- New row 94: "No lerp W" ---
$ws.Range("A94").Value = "No lerp W"
$ws.Range("B94").Value = 274
$ws.Range("C94").Formula = "=B94/30"
$ws.Range("D94").Formula = "=B94/B`$85"
$ws.Range("E94").Value = 252
$ws.Range("F94").Value = 909
$ws.Range("G94").Formula = "=E94*C94"
$ws.Range("H94").Formula = "=F94*C94"

# --- Row 96 (formerly row 94, "B2F") keeps its own independent formulas,
#     untouched values, just re-pointed at its new row number ---
$ws.Range("C96").Formula = "=B96/30"
$ws.Range("G96").Formula = "=E96*C96"
$ws.Range("H96").Formula = "=F96*C96"

# Update the saved selection to match.
$ws.Range("A94").Select()

$wb.Save()
